{"js": "// The commit merges each \"<id>...</id>\" marker (currently split across\n// three runs: the \"<id>\" open tag, the bare id text, and the \"</id>\"\n// close tag) into a single run per occurrence, e.g.\n//   <id>p068r_3</id>   (one run, formatted like the old \"<id>\" run)\n// instead of three separate runs. There are three such occurrences in\n// this document.\nconst body = context.document.body;\n\nconst openTags = body.search(\"<id>\", { matchCase: true, matchWildcards: false });\nconst closeTags = body.search(\"</id>\", { matchCase: true, matchWildcards: false });\nopenTags.load(\"items\");\ncloseTags.load(\"items\");\nawait context.sync();\n\nconst count = Math.min(openTags.items.length, closeTags.items.length);\nfor (let i = 0; i < count; i++) {\n  const openRange = openTags.items[i];\n  const closeRange = closeTags.items[i];\n\n  // Range spanning from the start of \"<id>\" to the end of \"</id>\",\n  // covering the id text in between regardless of how many runs it\n  // currently spans.\n  const fullRange = openRange.expandTo(closeRange);\n  fullRange.load(\"text\");\n  await context.sync();\n\n  const combinedText = fullRange.text;\n\n  // Re-inserting the same text over the whole range collapses the\n  // three runs into one run, which inherits the formatting of the\n  // range's leading run (i.e. the \"<id>\" run's Courier New / brown\n  // formatting) -- matching the target edit.\n  fullRange.insertText(combinedText, \"Replace\");\n  await context.sync();\n}\n", "ps1": "# The commit merges each \"<id>...</id>\" marker (currently split across\n# three runs: the \"<id>\" open tag, the bare id text, and the \"</id>\"\n# close tag) into a single run per occurrence, e.g.\n#   <id>p068r_3</id>   (one run, formatted like the old \"<id>\" run)\n# instead of three separate runs. There are three such occurrences in\n# this document.\n\n$d = $word.ActiveDocument\n$wdFindStop = 0\n\n$searchFrom = 0\n$docEnd = $d.Content.End\n\nwhile ($true) {\n    # Find the next \"<id>\" starting from $searchFrom.\n    $openRange = $d.Range($searchFrom, $docEnd)\n    $openRange.Find.ClearFormatting()\n    $openRange.Find.MatchCase = $true\n    $openRange.Find.MatchWildcards = $false\n    $openRange.Find.Wrap = $wdFindStop\n    $foundOpen = $openRange.Find.Execute(\"<id>\")\n    if (-not $foundOpen) { break }\n\n    $idStart = $openRange.Start\n\n    # Find the matching \"</id>\" after it.\n    $closeRange = $d.Range($openRange.End, $docEnd)\n    $closeRange.Find.ClearFormatting()\n    $closeRange.Find.MatchCase = $true\n    $closeRange.Find.MatchWildcards = $false\n    $closeRange.Find.Wrap = $wdFindStop\n    $foundClose = $closeRange.Find.Execute(\"</id>\")\n    if (-not $foundClose) { break }\n\n    # Range spanning from the start of \"<id>\" to the end of \"</id>\".\n    $full = $d.Range($idStart, $closeRange.End)\n    $combinedText = $full.Text\n\n    # Assigning a different placeholder first, then restoring the\n    # original combined text, forces Word to collapse the range's\n    # (possibly many) runs into a single run that inherits the\n    # formatting of the range's first run -- i.e. the old \"<id>\"\n    # run's Courier New / brown formatting -- matching the target edit.\n    $full.Text = \"@@@\"\n    $full.Text = $combinedText\n\n    $searchFrom = $full.End\n}\n"}
